# faturamento_diario.xlsx update
#
# The sheet holds a rolling window of daily revenue rows grouped by month
# (most recent month first). This commit:
#   1. Drops the oldest month in the window (02/2025 - rows 92:119).
#   2. Adds the two missing closing days of 05/2025 (day 30 and day 31)
#      right after the existing 05/2025 block.
#   3. Adds a new, most-recent month (06/2025) with a single day of data,
#      inserted at the very top of the data (row 2).
#
# All the other existing daily rows are left exactly as they were; they
# simply get shifted down as rows are inserted above them / removed below
# them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Remove the oldest month (02/2025), rows 92:119 -----------------
$ws.Rows("92:119").Delete()

# --- 2) Insert 05/2025 days 30 and 31 right after day 29 (row 31) ------
$ws.Rows("31:32").Insert()
$ws.Rows("31:32").ClearFormats()

$ws.Range("A31").Value2 = 30
$ws.Range("B31").Value2 = 25298.12
$ws.Range("C31").Value2 = 5
$ws.Range("D31").Value2 = 2025
$ws.Range("E31").Value2 = "05/2025"

$ws.Range("A32").Value2 = 31
$ws.Range("B32").Value2 = 19327.9
$ws.Range("C32").Value2 = 5
$ws.Range("D32").Value2 = 2025
$ws.Range("E32").Value2 = "05/2025"

# --- 3) Insert the new 06/2025 month (single day) at the very top ------
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").ClearFormats()

$ws.Range("A2").Value2 = 1
$ws.Range("B2").Value2 = 5140.01
$ws.Range("C2").Value2 = 6
$ws.Range("D2").Value2 = 2025
$ws.Range("E2").Value2 = "06/2025"
